$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("L2").Value = 13391.57
$ws.Range("M2").Value = 279350.01
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 63774.22
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 13391.57
$ws.Range("U2").Value = 279350.01
$ws.Range("V2").Value = 0
$ws.Range("Z2").Value = 63774.22
$ws.Range("AA2").Value = 0
$ws.Range("AC2").Value = 356515.8

# --- Row 4 ---
$ws.Range("L4").Value = 33666.72
$ws.Range("M4").Value = 341703.8
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 33666.72
$ws.Range("U4").Value = 341703.8
$ws.Range("V4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AC4").Value = 375370.52

# --- Row 14 ---
$ws.Range("M14").Value = 334593.33
$ws.Range("U14").Value = 334593.33
$ws.Range("AC14").Value = 502258.49

# --- Row 43 ---
$ws.Range("K43").Value = 74021.34
$ws.Range("L43").Value = 26512.21
$ws.Range("M43").Value = 337917.65
$ws.Range("N43").Value = 60888.94
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 0
$ws.Range("S43").Value = 74021.34
$ws.Range("T43").Value = 26512.21
$ws.Range("U43").Value = 337917.65
$ws.Range("V43").Value = 60888.94
$ws.Range("W43").Value = 0
$ws.Range("AA43").Value = 74021.34
$ws.Range("AC43").Value = 499340.14

# --- Row 47 ---
$ws.Range("L47").Value = 136868.28
$ws.Range("T47").Value = 136868.28
$ws.Range("AC47").Value = 722555.26
